# Daily attendance processing - reorder "Recorded By" (column G) entries
# so that "System" appears first in the comma-separated list whenever it
# was previously listed last among multiple recorders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $last = $parts[$parts.Length - 1]
    if ($last -eq "System") {
        $rest = $parts[0..($parts.Length - 2)]
        $newParts = @("System") + $rest
        $newVal = [string]::Join(", ", $newParts)
        $cell.Value = $newVal
    }
}
